$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.494.46'
$ws.Range('E2').Value = '  -2.92%  '
$ws.Range('D3').Value = '2.477.84'
$ws.Range('E3').Value = '  -2.09%  '
$ws.Range('E4').Value = '  +0.63%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.48'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.35'
$ws.Range('E6').Value = '  -6.75%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.544'
$ws.Range('E7').Value = '  -3.45%  '
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.493'
$ws.Range('E9').Value = '  -4.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.89'
$ws.Range('E10').Value = '  -6.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0779'
$ws.Range('E11').Value = '  -2.86%  '
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = '2.851.23'
$ws.Range('E13').Value = '  -2.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.85'
$ws.Range('E14').Value = '  -5.33%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.34'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.445.39'
$ws.Range('E16').Value = '  -2.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.783'
$ws.Range('E17').Value = '  -3.63%  '
$ws.Range('D18').Value = '41.383.01'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.28'
$ws.Range('E19').Value = '  -4.84%  '
$ws.Range('D20').Value = '0.0₃0920'
$ws.Range('E20').Value = '  -2.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.27'
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.01'
$ws.Range('E22').Value = '  -9.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.04'
$ws.Range('E23').Value = '  -2.99%  '
$ws.Range('E24').Value = '  -4.40%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.87'
$ws.Range('E26').Value = '  -6.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.00'
$ws.Range('E27').Value = '  -6.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.24'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.72'
$ws.Range('E29').Value = '  -2.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.16'
$ws.Range('E30').Value = '  -4.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '152.04'
$ws.Range('E31').Value = '  -2.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.43'
$ws.Range('E32').Value = '  -8.32%  '
$ws.Range('B33').Value = 'ApeXProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('E33').Value = '  -6.03%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('E34').Value = '  -3.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0750'
$ws.Range('E35').Value = '  -4.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.57'
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.99'
$ws.Range('E37').Value = '  -4.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.85'
$ws.Range('E38').Value = '  -6.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.113'
$ws.Range('E39').Value = '  -3.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0994'
$ws.Range('E40').Value = '  -8.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.03'
$ws.Range('E41').Value = '  -5.04%  '
$ws.Range('E42').Value = '  +0.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.40'
$ws.Range('E43').Value = '  -13.24%  '
$ws.Range('D44').Value = '1.960.55'
$ws.Range('E44').Value = '  -3.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0282'
$ws.Range('E45').Value = '  -5.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.95'
$ws.Range('E46').Value = '  -8.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.69'
$ws.Range('E47').Value = '  -2.07%  '
$ws.Range('D48').Value = '2.711.30'
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '68.40'
$ws.Range('E49').Value = '  -5.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '95.63'
$ws.Range('E50').Value = '  -4.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.176'
$ws.Range('E51').Value = '  -6.84%  '
